$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the confidentiality footer date string (2021-03-17 -> 2021-03-18)
$ws.Range("A80").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-18 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) columns for rows 2-77
$ws.Range("D2").Value = 0.0652076995478746
$ws.Range("E2").Value = -0.03390509778775253
$ws.Range("D3").Value = 0.03933448022505742
$ws.Range("E3").Value = -0.03435882553663749
$ws.Range("D4").Value = 0.03072534954316602
$ws.Range("E4").Value = -0.02666216672291599
$ws.Range("D5").Value = 0.02984058203920127
$ws.Range("E5").Value = -0.01058201058201069
$ws.Range("D6").Value = 0.02723620449893906
$ws.Range("E6").Value = 0.01650654458701406
$ws.Range("D7").Value = 0.02611929005820624
$ws.Range("E7").Value = -0.02923802480045345
$ws.Range("D8").Value = 0.02554477656555617
$ws.Range("E8").Value = -0.001866019779809713
$ws.Range("D9").Value = 0.1701254037514828
$ws.Range("E9").Value = -0.01074053137365738
$ws.Range("D10").Value = 0.02252460847428141
$ws.Range("E10").Value = -0.01577884372633176
$ws.Range("D11").Value = 0.02255253969794154
$ws.Range("E11").Value = 0.002180345740538936
$ws.Range("D12").Value = 0.02122538844199442
$ws.Range("E12").Value = -0.01516082769924199
$ws.Range("D13").Value = 0.01999382218198602
$ws.Range("E13").Value = 0.02608695652173898
$ws.Range("D14").Value = 0.01737439188645787
$ws.Range("E14").Value = -0.03615675628844539
$ws.Range("D15").Value = 0.01610243406648586
$ws.Range("E15").Value = -0.040051518550833
$ws.Range("D16").Value = 0.01715119297643123
$ws.Range("E16").Value = -0.01254552812626486
$ws.Range("D17").Value = 0.01548694362894517
$ws.Range("E17").Value = -0.01716056849107961
$ws.Range("D18").Value = 0.01452641059014185
$ws.Range("E18").Value = -0.008157461457865445
$ws.Range("D19").Value = 0.01365347622434449
$ws.Range("E19").Value = -0.04311942058278584
$ws.Range("D20").Value = 0.01185856059850825
$ws.Range("E20").Value = -0.01209063214013728
$ws.Range("D21").Value = 0.01306290645941074
$ws.Range("E21").Value = -0.0189782049927818
$ws.Range("D22").Value = 0.01266706080912557
$ws.Range("E22").Value = -0.03084398436674762
$ws.Range("D23").Value = 0.01190597677759597
$ws.Range("E23").Value = -0.03103862443896588
$ws.Range("D24").Value = 0.01285881618593003
$ws.Range("E24").Value = -0.004643449419568735
$ws.Range("D25").Value = 0.01178062077679794
$ws.Range("E25").Value = 0.02802544153557851
$ws.Range("D26").Value = 0.01109421513667103
$ws.Range("E26").Value = -0.05213170113972143
$ws.Range("D27").Value = 0.01068745623881805
$ws.Range("E27").Value = -0.02143974960876371
$ws.Range("D28").Value = 0.01072684428176393
$ws.Range("E28").Value = -0.031240498631803
$ws.Range("D29").Value = 0.0105787418952307
$ws.Range("E29").Value = -0.02834782608695663
$ws.Range("D30").Value = 0.00924924909957554
$ws.Range("E30").Value = -0.02502667221207577
$ws.Range("D31").Value = 0.01036382200060032
$ws.Range("E31").Value = 0.002178649237472685
$ws.Range("D32").Value = 0.01053805764280358
$ws.Range("E32").Value = 0.0111773738528016
$ws.Range("D33").Value = 0.009253472233691818
$ws.Range("E33").Value = 0.007094281647499923
$ws.Range("D34").Value = 0.00955745426364568
$ws.Range("E34").Value = 0.0100448430493274
$ws.Range("D35").Value = 0.008803478477658328
$ws.Range("E35").Value = -0.06932075633006074
$ws.Range("D36").Value = 0.009212787981264705
$ws.Range("E36").Value = -0.01307572209211549
$ws.Range("D37").Value = 0.009194724674945576
$ws.Range("E37").Value = 0.008412914961346107
$ws.Range("D38").Value = 0.008655334277916026
$ws.Range("E38").Value = -0.007536231884057942
$ws.Range("D39").Value = 0.00808817991006263
$ws.Range("E39").Value = -0.04019934241816414
$ws.Range("D40").Value = 0.007868911441688758
$ws.Range("E40").Value = 0.004165958170379236
$ws.Range("D41").Value = 0.007839474942502029
$ws.Range("E41").Value = -0.03113799283154117
$ws.Range("D42").Value = 0.008477544512941635
$ws.Range("E42").Value = -0.01139345394282565
$ws.Range("D43").Value = 0.008403618759302236
$ws.Range("E43").Value = -0.0003880983182406039
$ws.Range("D44").Value = 0.007752837973304724
$ws.Range("E44").Value = 0.007852612503775225
$ws.Range("D45").Value = 0.008130745756897985
$ws.Range("E45").Value = -0.03619364991206142
$ws.Range("D46").Value = 0.00808107166452038
$ws.Range("E46").Value = -0.0005588153115393402
$ws.Range("D47").Value = 0.007743137308800006
$ws.Range("E47").Value = 0.0111456713322966
$ws.Range("D48").Value = 0.006246224424020327
$ws.Range("E48").Value = -0.01499491244042195
$ws.Range("D49").Value = 0.007210855157312708
$ws.Range("E49").Value = -0.001333688983729209
$ws.Range("D50").Value = 0.006928030611149306
$ws.Range("E50").Value = -0.02571066449393422
$ws.Range("D51").Value = 0.006838174024853452
$ws.Range("E51").Value = -0.01262680306467479
$ws.Range("D52").Value = 0.006498650767188341
$ws.Range("E52").Value = 0.00635049317659786
$ws.Range("D53").Value = 0.00632717479678383
$ws.Range("E53").Value = -0.02498017446471035
$ws.Range("D54").Value = 0.00608984302209083
$ws.Range("E54").Value = -0.006165719150806104
$ws.Range("D55").Value = 0.005687282170428683
$ws.Range("E55").Value = -0.009788525501336465
$ws.Range("D56").Value = 0.005373833629940907
$ws.Range("E56").Value = -0.002801120448179373
$ws.Range("D57").Value = 0.005623709367355526
$ws.Range("E57").Value = 0.0596002855103499
$ws.Range("D58").Value = 0.005539413937866257
$ws.Range("E58").Value = 0.003260869565217428
$ws.Range("D59").Value = 0.005000525299345572
$ws.Range("E59").Value = -0.00570272259013993
$ws.Range("D60").Value = 0.005455996585767501
$ws.Range("E60").Value = -0.0495459248189446
$ws.Range("D61").Value = 0.004704822284787977
$ws.Range("E61").Value = 0.001599715606114493
$ws.Range("D62").Value = 0.004802832446852882
$ws.Range("E62").Value = 0.005989692157682081
$ws.Range("D63").Value = 0.004433872690001041
$ws.Range("E63").Value = -0.001131648434553045
$ws.Range("D64").Value = 0.004159912544160916
$ws.Range("E64").Value = 0.002412351238340271
$ws.Range("D65").Value = 0.004003363889395131
$ws.Range("E65").Value = -0.0006684491978610207
$ws.Range("D66").Value = 0.003759425627668743
$ws.Range("E66").Value = -0.01955288621955276
$ws.Range("D67").Value = 0.003832431490708557
$ws.Range("E67").Value = 0.007200837915684932
$ws.Range("D68").Value = 0.003154430305605135
$ws.Range("E68").Value = -0.05193462440847807
$ws.Range("D69").Value = 0.003334687049914777
$ws.Range("E69").Value = -0.008764670478483283
$ws.Range("D70").Value = 0.003388876968872165
$ws.Range("E70").Value = -0.02334419109663421
$ws.Range("D71").Value = 0.002673369335231106
$ws.Range("E71").Value = -0.01707957957957973
$ws.Range("D72").Value = 0.00223136190213038
$ws.Range("E72").Value = -0.04637871263937032
$ws.Range("D73").Value = 0.002192851936575015
$ws.Range("E73").Value = -0.03746853786896509
$ws.Range("D74").Value = 0.001885817542358892
$ws.Range("E74").Value = -0.02623001707279216
$ws.Range("D75").Value = 0.001768615117329913
$ws.Range("E75").Value = -0.05848976310936693
$ws.Range("D76").Value = 0.001693518593836497
$ws.Range("E76").Value = -0.02281368821292773
$ws.Range("D77").Value = 1
$ws.Range("E77").Value = -0.01445436224958985

$ws.Protect()
